# Update the NATMI LR-pairs worksheet with new TPM-based values.
# This replaces the existing data rows (2-9) and appends new rows (10-13)
# for the Efnb2-Ephb1 ligand-receptor pair across all sending/target cluster
# combinations, reflecting the updated TPM calculations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $data = @(
        @("ECs", "Efnb2", "Ephb1", "ECs", 3, 1, 45.71598933333333, 137.147968, 0.6549002937372808, 0.6549002937372808, 3, 1, 4.002008666666667, 12.006026, 0.4834231243738785, 0.4834231243738787, 182.9557855172409, 1646.602069655168, 0.316593946151847, 0.3165939461518472),
        @("ECs", "Efnb2", "Ephb1", "MuSCs", 3, 1, 45.71598933333333, 137.147968, 0.6549002937372808, 0.6549002937372808, 3, 1, 4.265473333333333, 12.79642, 0.5152483708764571, 0.5152483708764573, 195.0003334082845, 1755.00300067456, 0.3374363094346471, 0.3374363094346473),
        @("ECs", "Efnb2", "Ephb1", "Resolving-Mac", 3, 1, 45.71598933333333, 137.147968, 0.6549002937372808, 0.6549002937372808, 1, 0.3333333333333333, 0.010998, 0.032994, 0.001328504749664189, 0.00132850474966419, 0.502784450688, 4.525060056192, 0.0008700381507864502, 0.0008700381507864506),
        @("FAPs", "Efnb2", "Ephb1", "ECs", 3, 1, 12.691493, 38.074479, 0.1818108415648851, 0.1818108415648851, 3, 1, 4.002008666666667, 12.006026, 0.4834231243738785, 0.4834231243738787, 50.79146497893933, 457.1231848104539, 0.08789156507434097, 0.08789156507434098),
        @("FAPs", "Efnb2", "Ephb1", "MuSCs", 3, 1, 12.691493, 38.074479, 0.1818108415648851, 0.1818108415648851, 3, 1, 4.265473333333333, 12.79642, 0.5152483708764571, 0.5152483708764573, 54.13522495168667, 487.2170245651799, 0.0936777399239847, 0.0936777399239847),
        @("FAPs", "Efnb2", "Ephb1", "Resolving-Mac", 3, 1, 12.691493, 38.074479, 0.1818108415648851, 0.1818108415648851, 1, 0.3333333333333333, 0.010998, 0.032994, 0.001328504749664189, 0.00132850474966419, 0.139581040014, 1.256229360126, 0.0002415365665593932, 0.0002415365665593933),
        @("MuSCs", "Efnb2", "Ephb1", "ECs", 3, 1, 11.24784666666667, 33.74354, 0.161130015850732, 0.161130015850732, 3, 1, 4.002008666666667, 12.006026, 0.4834231243738785, 0.4834231243738787, 45.01397984133778, 405.12581857204, 0.07789397569297342, 0.07789397569297345),
        @("MuSCs", "Efnb2", "Ephb1", "MuSCs", 3, 1, 11.24784666666667, 33.74354, 0.161130015850732, 0.161130015850732, 3, 1, 4.265473333333333, 12.79642, 0.5152483708764571, 0.5152483708764573, 47.97739001408889, 431.7965101267999, 0.08302197816638736, 0.08302197816638739),
        @("MuSCs", "Efnb2", "Ephb1", "Resolving-Mac", 3, 1, 11.24784666666667, 33.74354, 0.161130015850732, 0.161130015850732, 1, 0.3333333333333333, 0.010998, 0.032994, 0.001328504749664189, 0.00132850474966419, 0.12370381764, 1.11333435876, 0.0002140619913711635, 0.0002140619913711636),
        @("Resolving-Mac", "Efnb2", "Ephb1", "ECs", 2, 0.6666666666666666, 0.1507006666666667, 0.452102, 0.00215884884710222, 0.00215884884710222, 3, 1, 4.002008666666667, 12.006026, 0.4834231243738785, 0.4834231243738787, 0.6031053740724445, 5.427948366652, 0.001043637454717101, 0.001043637454717101),
        @("Resolving-Mac", "Efnb2", "Ephb1", "MuSCs", 2, 0.6666666666666666, 0.1507006666666667, 0.452102, 0.00215884884710222, 0.00215884884710222, 3, 1, 4.265473333333333, 12.79642, 0.5152483708764571, 0.5152483708764573, 0.6428096749822223, 5.785287074839999, 0.001112343351437937, 0.001112343351437937),
        @("Resolving-Mac", "Efnb2", "Ephb1", "Resolving-Mac", 2, 0.6666666666666666, 0.1507006666666667, 0.452102, 0.00215884884710222, 0.00215884884710222, 1, 0.3333333333333333, 0.010998, 0.032994, 0.001328504749664189, 0.00132850474966419, 0.001657405932, 0.014916653388, 0.000002868040947182358, 0.000002868040947182359)
    )
$startRow = 2
$r = $startRow
foreach ($rowValues in $data) {
    $c = 1
    foreach ($val in $rowValues) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
Write-Host "Updated Efnb2-Ephb1 LR-pairs data (rows 2-13)."
